$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.895.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.546.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.766.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.549.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.873.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "213.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  -3.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.11%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.361.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.805"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.680.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0505"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0963"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0947"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
